$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
# Force a real text-model change so the three split runs ("Below", " ",
# "section-level") collapse into a single run, matching the target XML.
$tr.Text = "_tmp_"
$tr.Text = "Below section-level"
